$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The row that used to be "State transition Diagram" (row 11) was cut and
# re-inserted above row 8, pushing "Pairwise Testing" / "Use Case Diagram" /
# "Experience Based Techniques" down by one row each, and the trailing blank
# bordered row slides from row 15 to row 16.
# ---------------------------------------------------------------------------

# Snapshot the current text of the four rows that move around (use .Text,
# not .Value, as the getter to read back to a plain PowerShell variable).
$moveUp   = $ws.Range("A11").Text
$row8Val  = $ws.Range("A8").Text
$row9Val  = $ws.Range("A9").Text
$row10Val = $ws.Range("A10").Text

# Shift rows 8-10 down into 9-11, then drop the moved value into row 8.
$ws.Range("A11").Value = $row10Val
$ws.Range("A10").Value = $row9Val
$ws.Range("A9").Value  = $row8Val
$ws.Range("A8").Value  = $moveUp

# The cursor/selection ends up on the destination row (row 8), selected as a
# whole row, matching an "insert cut cells" drop target.
$ws.Rows("8:8").Select()

# ---------------------------------------------------------------------------
# Move the trailing empty bordered cell from row 15 down to row 16.
# ---------------------------------------------------------------------------
$ws.Rows("15:15").Delete()
$ws.Rows("16:16").Insert()
$tail = $ws.Range("A16")
$tail.Font.Name = "Times New Roman"
$tail.Font.Size = 14
$tail.HorizontalAlignment = -4108
$tail.VerticalAlignment = -4108
$tail.Borders.LineStyle = 1
$tail.Borders.Weight = -4138
